# This workbook lists modules; the edit moves the "CodeModule" column (F)
# to the front of the table (column A), shifting the other columns one
# place to the right, and renumbers the module codes from the old
# "AP4x" scheme to the new "G3EI1x" scheme to match the new project name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowCount = 10
$colCount = 6

# Capture the existing table contents (A1:F10) before moving anything,
# cell by cell, so we do not disturb column/row formatting with a bulk
# cut/insert operation.
$data = @{}
for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $data["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# New codes replacing the old "AP41".."AP49" values (header stays "CodeModule").
$newCodes = @("CodeModule", "G3EI11", "G3EI12", "G3EI13", "G3EI14", "G3EI15", "G3EI16", "G3EI17", "G3EI18", "G3EI19")

# Write back the data shifted one column to the right, with the former
# column F (CodeModule) now occupying column A using the updated codes.
for ($r = 1; $r -le $rowCount; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $newCodes[$r - 1]
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $data["$r,$c"]
    }
}

# Update the active selection to reflect where the user left off editing.
$ws.Range("D7").Select()
